$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.375.15'
$ws.Range("E2").Value = '  +4.08%  '

$ws.Range("D3").Value = '2.045.12'
$ws.Range("E3").Value = '  +2.62%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.56%  '

$ws.Range("E6").Value = '  +1.36%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.48'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +9.67%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.395'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0832'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.12%  '

$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.916'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +23.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").Value = '2.341.07'
$ws.Range("E16").Value = '  +2.63%  '

$ws.Range("E17").Value = '  +5.52%  '

$ws.Range("D18").Value = '2.041.63'
$ws.Range("E18").Value = '  +2.63%  '

$ws.Range("D19").Value = '37.241.22'
$ws.Range("E19").Value = '  +3.82%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0915'
$ws.Range("E20").Value = '  +7.57%  '

$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.41%  '

$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +23.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.03%  '

$ws.Range("E33").Value = '  +4.28%  '

$ws.Range("E34").Value = '  +4.15%  '

$ws.Range("E35").Value = '  +4.51%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.37'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.38%  '

$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.01%  '

$ws.Range("E39").Value = '  +3.42%  '

$ws.Range("E40").Value = '  +30.39%  '

$ws.Range("E41").Value = '  +4.56%  '

$ws.Range("E42").Value = '  +6.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.29%  '

$ws.Range("E44").Value = '  +4.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.15%  '

$ws.Range("E46").Value = '  +1.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("D49").Value = '1.389.86'
$ws.Range("E49").Value = '  +2.01%  '

$ws.Range("E50").Value = '  +1.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.14%  '
